$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right above the current row 738 (shifts the
# existing 738:849 block down to 740:851, carrying their values/format with
# them, and also pushes the sheet dimension out to R851 automatically).
$ws.Rows("738:739").Insert()

# --- New row 738: Brócoli, Primera, Región Metropolitana, fecha 44776 ---
$ws.Range("A738").Value = 6
$ws.Range("B738").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C738").Value = "Metropolitana"
$ws.Range("D738").Value = 44776
$ws.Range("E738").Value = 13
$ws.Range("F738").Value = 100112023
$ws.Range("G738").Value = "Brócoli"
$ws.Range("H738").Value = "Sin especificar"
$ws.Range("I738").Value = "Primera"
$ws.Range("J738").Value = 12100
$ws.Range("K738").Value = 500
$ws.Range("L738").Value = 600
$ws.Range("M738").Value = 549
$ws.Range("N738").Value = "$/unidad"
$ws.Range("O738").Value = "Región Metropolitana"
$ws.Range("P738").Value = 549
$ws.Range("Q738").Value = 1
$ws.Range("R738").Value = "Hortaliza"

# --- New row 739: Brócoli, Segunda, Región Metropolitana, fecha 44776 ---
$ws.Range("A739").Value = 6
$ws.Range("B739").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C739").Value = "Metropolitana"
$ws.Range("D739").Value = 44776
$ws.Range("E739").Value = 13
$ws.Range("F739").Value = 100112023
$ws.Range("G739").Value = "Brócoli"
$ws.Range("H739").Value = "Sin especificar"
$ws.Range("I739").Value = "Segunda"
$ws.Range("J739").Value = 4700
$ws.Range("K739").Value = 300
$ws.Range("L739").Value = 400
$ws.Range("M739").Value = 355
$ws.Range("N739").Value = "$/unidad"
$ws.Range("O739").Value = "Región Metropolitana"
$ws.Range("P739").Value = 355
$ws.Range("Q739").Value = 1
$ws.Range("R739").Value = "Hortaliza"
